$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- New rows 20-26: intro/stats/methods courses --
$ws.Range("A20").Value = "Intro 1"
$ws.Range("B20").Value = "Cog Psych"
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Formula = "=(D20*C20)"

$ws.Range("A21").Value = "Intro 2"
$ws.Range("B21").Value = "Dev Pysch"
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 3
$ws.Range("E21").Formula = "=(D21*C21)"

$ws.Range("A22").Value = "Intro 3"
$ws.Range("B22").Value = "Social Psych"
$ws.Range("C22").Value = 3.7
$ws.Range("D22").Value = 3
$ws.Range("E22").Formula = "=(D22*C22)"

$ws.Range("A23").Value = "Stats I"
$ws.Range("B23").Value = "Stats I"
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 4
$ws.Range("E23").Formula = "=(D23*C23)"

$ws.Range("A24").Value = "Stats II"
$ws.Range("B24").Value = "Stats II"
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 4
$ws.Range("E24").Formula = "=(D24*C24)"

$ws.Range("A25").Value = "Calc I"
$ws.Range("B25").Value = "Calc I"
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Formula = "=(D25*C25)"

$ws.Range("A26").Value = "Research Methods"
$ws.Range("B26").Value = "RM"
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 3
$ws.Range("E26").Formula = "=(D26*C26)"

# -- 200/300-level courses (A holds plain numeric level, not text) --
$ws.Range("A27").Value = 200
$ws.Range("B27").Value = "Abnormal"
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 3
$ws.Range("E27").Formula = "=(D27*C27)"

$ws.Range("A28").Value = 200
$ws.Range("B28").Value = "Positive"
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 3
$ws.Range("E28").Formula = "=(D28*C28)"

$ws.Range("A29").Value = 300
$ws.Range("B29").Value = "Foundations"
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 4
$ws.Range("E29").Formula = "=(D29*C29)"

$ws.Range("A30").Value = 300
$ws.Range("B30").Value = "Evolutionary"
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 3
$ws.Range("E30").Formula = "=(D30*C30)"

$ws.Range("A31").Value = 300
$ws.Range("B31").Value = "Theory"
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 3
$ws.Range("E31").Formula = "=(D31*C31)"

# -- "3 N" (300-level numbered) courses --
$ws.Range("A32").Value = "3 N"
$ws.Range("B32").Value = "Planets"
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = 3
$ws.Range("E32").Formula = "=(D32*C32)"

$ws.Range("A33").Value = "3 N"
$ws.Range("B33").Value = "Cog Neuro"
$ws.Range("C33").Value = 4
$ws.Range("D33").Value = 3
$ws.Range("E33").Formula = "=(D33*C33)"

$ws.Range("A34").Value = "3 N"
$ws.Range("B34").Value = "Visual Cognition"
$ws.Range("C34").Value = 4
$ws.Range("D34").Value = 3
$ws.Range("E34").Formula = "=(D34*C34)"

# -- Overall GPA summary for the new block --
$ws.Range("D36").Formula = "=(SUM(E20:E34))/(SUM(D20:D34))"

# -- View state: scroll down a bit and move the selection past the new data --
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("D37").Select()
